$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 13, shifting existing rows 13:57 down to 14:58.
$ws.Rows.Item(13).Insert(-4121)

# Populate the newly inserted row 13 with the new weekly record.
# Columns A,B,C,E,F,G,H,I,J,K repeat the same constant values used throughout the table.
$ws.Range("A13").Value = 8
$ws.Range("B13").Value = "Terminal La Palmera de La Serena"
$ws.Range("C13").Value = "Coquimbo"
$ws.Range("D13").Value = 44707
$ws.Range("E13").Value = 4
$ws.Range("F13").Value = "Fruta"
$ws.Range("G13").Value = 100104
$ws.Range("H13").Value = "Frutos de pepita"
$ws.Range("I13").Value = 100104003
$ws.Range("J13").Value = "Membrillo"
$ws.Range("K13").Value = "Champion"
$ws.Range("L13").Value = "Primera"
$ws.Range("M13").Value = 20
$ws.Range("N13").Value = 280000
$ws.Range("O13").Value = 290000
$ws.Range("P13").Value = 285000
$ws.Range("Q13").Value = "$/bins (450 kilos)"
$ws.Range("R13").Value = "Región de O'Higgins"
$ws.Range("S13").Value = 633
$ws.Range("T13").Value = 450

# Give the date cell the same number format used by the rest of column D.
$ws.Range("D13").NumberFormat = $ws.Range("D14").NumberFormat()
